$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1004.5
$ws.Range("I15").Value = 1004.5
$ws.Range("K15").Value = 3013.5
$ws.Range("M15").Value = -2844.5
$ws.Range("H33").Value = 107.666664
$ws.Range("H62").Value = 1521.8
$ws.Range("I62").Value = 1521.8
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1521.8
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -897.8
$ws.Range("H65").Value = 1521.8
$ws.Range("I65").Value = 1521.8
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7609
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -4489
$ws.Range("H69").Value = 5266.3335
$ws.Range("I69").Value = 5266.3335
$ws.Range("K69").Value = 15799.0005
$ws.Range("M69").Value = -14925.0005
$ws.Range("H72").Value = 5266.3335
$ws.Range("I72").Value = 5266.3335
$ws.Range("K72").Value = 47397.0015
$ws.Range("M72").Value = -43029.0015
$ws.Range("H98").Value = 3783.5
$ws.Range("I98").Value = 2914.0833
$ws.Range("J98").Value = 9000
$ws.Range("K98").Value = 2914.0833
$ws.Range("L98").Value = 9000
$ws.Range("M98").Value = -1416.0833
$ws.Range("N98").Value = -11996
$ws.Range("H122").Value = 3783.5
$ws.Range("I122").Value = 2914.0833
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 8742.249899999999
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -6292.249899999999
$ws.Range("N122").Value = -31900
$ws.Range("H127").Value = 6244.3335
$ws.Range("I127").Value = 6244.3335
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 18733.0005
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -13773.0005
$ws.Range("H132").Value = 1379.6875
$ws.Range("I132").Value = 1379.6875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4139.0625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1609.0625
$ws.Range("H137").Value = 1771.7
$ws.Range("I137").Value = 1237.625
$ws.Range("J137").Value = 2127.75
$ws.Range("K137").Value = 3712.875
$ws.Range("L137").Value = 6383.25
$ws.Range("M137").Value = -1162.875
$ws.Range("N137").Value = -11483.25
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4217.75
$ws.Range("I32").Value = 2776.389
$ws.Range("K32").Value = 2776.389
$ws.Range("M32").Value = -2489.389
$ws.Range("H45").Value = 2158.6365
$ws.Range("I45").Value = 1025
$ws.Range("J45").Value = 3519
$ws.Range("K45").Value = 1025
$ws.Range("L45").Value = 3519
$ws.Range("M45").Value = -648
$ws.Range("N45").Value = -4273
$ws.Range("H97").Value = 423.33334
$ws.Range("I97").Value = 423.33334
$ws.Range("K97").Value = 423.33334
$ws.Range("M97").Value = 72.66665999999998
$ws.Range("H132").Value = 1511.7273
$ws.Range("J132").Value = 4999
$ws.Range("L132").Value = 14997
$ws.Range("N132").Value = -20057

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2702.7778
$ws.Range("J107").Value = 3658.3333
$ws.Range("L107").Value = 3658.3333
$ws.Range("N107").Value = -7498.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3140.4092
$ws.Range("I31").Value = 1379.4
$ws.Range("J31").Value = 6914
$ws.Range("K31").Value = 1379.4
$ws.Range("L31").Value = 6914
$ws.Range("M31").Value = -1084.4
$ws.Range("N31").Value = -7504
$ws.Range("H34").Value = 3140.4092
$ws.Range("I34").Value = 1379.4
$ws.Range("J34").Value = 6914
$ws.Range("K34").Value = 1379.4
$ws.Range("L34").Value = 6914
$ws.Range("M34").Value = -1177.4
$ws.Range("N34").Value = -7318
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H140").Value = 59000
$ws.Range("J140").Value = 59000
$ws.Range("L140").Value = 59000
$ws.Range("N140").Value = -69360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7771.89
$ws.Range("I131").Value = 829.5
$ws.Range("J131").Value = 7913.5713
$ws.Range("K131").Value = 2488.5
$ws.Range("L131").Value = 23740.7139
$ws.Range("M131").Value = 2551.5
$ws.Range("N131").Value = -33820.7139
$ws.Range("H134").Value = 1848.6818
$ws.Range("I134").Value = 1426.5555
$ws.Range("K134").Value = 4279.666499999999
$ws.Range("M134").Value = 790.3335000000006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 12552505
$ws.Range("H30").Value = 12552505
$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 614.2857
$ws.Range("K97").Value = 614.2857
$ws.Range("M97").Value = -118.2857
$ws.Range("H132").Value = 3875.8809
$ws.Range("I132").Value = 3165
$ws.Range("J132").Value = 5653.0835
$ws.Range("K132").Value = 9495
$ws.Range("L132").Value = 16959.2505
$ws.Range("M132").Value = -6965
$ws.Range("N132").Value = -22019.2505

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("H55").Value = 242
$ws.Range("I55").Value = 214
$ws.Range("K55").Value = 214
$ws.Range("M55").Value = -41
$ws.Range("H132").Value = 2069.5908
$ws.Range("I132").Value = 1542.5714
$ws.Range("K132").Value = 4627.7142
$ws.Range("M132").Value = -2097.7142
$ws.Range("H133").Value = 87663
$ws.Range("J133").Value = 87663
$ws.Range("L133").Value = 87663
$ws.Range("N133").Value = -92723
$ws.Range("N29").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3286
$ws.Range("J132").Value = 3783.6
$ws.Range("L132").Value = 11350.8
$ws.Range("N132").Value = -16410.8
$ws.Range("H137").Value = 59000
$ws.Range("J137").Value = 59000
$ws.Range("L137").Value = 59000
$ws.Range("N137").Value = -69200
